# Add a new Job Posting row (Job_Id = JD_007) to the LinkedIn job posting sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the repeated Job_Description text, matching the pattern already used
# elsewhere in this sheet (the paragraph pasted multiple times in a row).
$paragraph = "We are seeking a Software Engineer to build and maintain high-quality software solutions." + "`n" + `
    "Work with global teams to drive innovation and deliver scalable applications." + "`n" + `
    "Join Akkodis and be part of a tech-driven, collaborative environment."

$jobDescription = ""
for ($i = 0; $i -lt 27; $i++) {
    $jobDescription += $paragraph
}

# New row goes right after the existing last row (row 7), i.e. row 8.
$newRow = 8

$ws.Cells.Item($newRow, 1).Value = "JD_007"
$ws.Cells.Item($newRow, 2).Value = "Senior Y Engineer"
$ws.Cells.Item($newRow, 3).Value = $jobDescription
$ws.Cells.Item($newRow, 4).Value = 1
$ws.Cells.Item($newRow, 5).Value = 5

# Keep the row height at the sheet default (avoid auto row-height growth
# from the long, multi-line Job_Description text).
$ws.Rows.Item($newRow).RowHeight = 15
